$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains numeric-looking values that must stay as literal text
# (to preserve exact formatting / trailing zeros), so a leading apostrophe
# is used to force Excel to store them as text, same as typing them by hand.

$ws.Range("D2").Value = "'248.02"
$ws.Range("D3").Value = "'21.62"
$ws.Range("D4").Value = "'5.520"
$ws.Range("D5").Value = "'0.05695"
$ws.Range("D6").Value = "'3.386"
$ws.Range("D7").Value = "'0.8071"
$ws.Range("D9").Value = "'0.1501"
$ws.Range("D10").Value = "'0.08289"
$ws.Range("D11").Value = "'0.03142"
$ws.Range("D12").Value = "'0.03015"
$ws.Range("D13").Value = "'0.09293"
$ws.Range("D14").Value = "'3.481"
$ws.Range("D15").Value = "'0.001660"
$ws.Range("D16").Value = "'0.04693"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005865"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006354"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005041"
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001043"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "UpBots"
$ws.Range("C22").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D22").Value = "'0.0003203"
$ws.Range("E22").Value = "21UpBotsUBXT"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.767"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "KuCoinToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D24").Value = "'6.428"
$ws.Range("E24").Value = "23KuCoinTokenKCS"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.120"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("D26").Value = "'0.3316"
$ws.Range("D40").Value = "'0.04119"
$ws.Range("D41").Value = "'0.006981"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1045"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002971"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.009147"
$ws.Range("D45").Value = "'0.00005887"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.0005504"
$ws.Range("D48").Value = "'0.6830"
$ws.Range("D49").Value = "'0.008646"
$ws.Range("D50").Value = "'0.00002102"
